# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# Update the "Date" column (BF) from "6-4-2011-12" to the correct "2012-06-04".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    if ($cell.Text -eq "6-4-2011-12") {
        # Assign via a literal-text formula, then collapse it back to a
        # plain value (copy / paste-values). This avoids Excel's automatic
        # date-recognition that a direct string assignment of "2012-06-04"
        # would trigger (which would reformat the cell as a date serial).
        $cell.Value = '="2012-06-04"'
        $cell.Copy()
        $cell.PasteSpecial(-4163)
    }
}
